# eventsliders.xlsx edit:
#  - adds labeled alarm sets
#  - adds Ramp/Soak pattern labels
#  - adds option to load Ramp/Soak patterns from background profile
#  - displays path the Ramp/Soak patterns were loaded from

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

# 1) pidRS(<int>) -> pidRS(<rs>) ; update its description to mention labels too
$ws.Range("B76").Value = "pidRS(<rs>)"
$ws.Range("C76").Value = "activates the PID Ramp-Soak pattern number <rs> (1-based!) or the one labeled <rs>"

# 2) palette(<int>) -> palette(<p>) ; update its description
$ws.Range("B83").Value = "palette(<p>)"
$ws.Range("C83").Value = "activates palette <p> with <p> either a number 0-9 or a palette label"

# 3) Insert a new row right after "clearBackground" (row 87) for the new
#    "alarmset(<as>)" command, pushing the "RC Command" section (and
#    everything below it) down by one row.
$ws.Rows.Item(88).Insert()
$ws.Range("B88").Value = "alarmset(<as>)"
$ws.Range("C88").Value = "activates the alarmset with the given number or label"
$ws.Rows.Item(88).RowHeight = 13.8

# Leave the view/selection pointing at the (now shifted) palette row, and
# keep "Commands" the active sheet/tab, mirroring the saved cursor state.
$ws1 = $wb.Worksheets.Item("Sliders")
$ws1.Range("B6").Select() | Out-Null
$ws.Range("B83:C83").Select() | Out-Null
$ws.Activate() | Out-Null
